$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for rows 2-8 from 2023-09-16 (45185) to 2023-10-05 (45204)
$ws.Range("C2:C8").Value = 45204
